# Fruta / hortaliza, semanal
# Insert two new weekly data rows (Navel Late, fecha 2021-09-13) right after
# the current row 129, pushing the existing rows 130:222 down to 132:224.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows at position 130 (shifts everything below it down by 2).
$ws.Rows("130:131").Insert()

# New row 130
$ws.Range("A130").Value2 = 4
$ws.Range("B130").Value2 = "Feria Lagunitas de Puerto Montt"
$ws.Range("C130").Value2 = "Los Lagos"
$ws.Range("D130").Value2 = 44452
$ws.Range("E130").Value2 = 10
$ws.Range("F130").Value2 = "Fruta"
$ws.Range("G130").Value2 = 100102
$ws.Range("H130").Value2 = "Cítricos"
$ws.Range("I130").Value2 = 100102005
$ws.Range("J130").Value2 = "Naranja"
$ws.Range("K130").Value2 = "Navel Late"
$ws.Range("L130").Value2 = "Primera"
$ws.Range("M130").Value2 = 200
$ws.Range("N130").Value2 = 13000
$ws.Range("O130").Value2 = 13000
$ws.Range("P130").Value2 = 13000
$ws.Range("Q130").Value2 = "`$/caja 15 kilos empedrada"
$ws.Range("R130").Value2 = "Región de O'Higgins"
$ws.Range("S130").Value2 = 867
$ws.Range("T130").Value2 = 15

# New row 131
$ws.Range("A131").Value2 = 4
$ws.Range("B131").Value2 = "Feria Lagunitas de Puerto Montt"
$ws.Range("C131").Value2 = "Los Lagos"
$ws.Range("D131").Value2 = 44452
$ws.Range("E131").Value2 = 10
$ws.Range("F131").Value2 = "Fruta"
$ws.Range("G131").Value2 = 100102
$ws.Range("H131").Value2 = "Cítricos"
$ws.Range("I131").Value2 = 100102005
$ws.Range("J131").Value2 = "Naranja"
$ws.Range("K131").Value2 = "Navel Late"
$ws.Range("L131").Value2 = "Segunda"
$ws.Range("M131").Value2 = 100
$ws.Range("N131").Value2 = 10000
$ws.Range("O131").Value2 = 10000
$ws.Range("P131").Value2 = 10000
$ws.Range("Q131").Value2 = "`$/caja 15 kilos empedrada"
$ws.Range("R131").Value2 = "Región de O'Higgins"
$ws.Range("S131").Value2 = 667
$ws.Range("T131").Value2 = 15
